$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 6 (hunk 1)
$ws.Range("H6").Value = 9889.666999999999
$ws.Range("I6").Value = 14434.333
$ws.Range("K6").Value = 43302.999
$ws.Range("M6").Value = -43190.999
# row 15 (hunk 2)
$ws.Range("H15").Value = 1919.4944
$ws.Range("I15").Value = 1919.4944
$ws.Range("K15").Value = 5758.483200000001
$ws.Range("M15").Value = -5589.483200000001
# row 51 (hunk 3)
$ws.Range("H51").Value = 61497.25
$ws.Range("J51").Value = 61497.25
$ws.Range("L51").Value = 61497.25
$ws.Range("N51").Value = -62465.25
# row 70 (hunk 4)
$ws.Range("H70").Value = 294670.6
$ws.Range("I70").Value = 333768.34
$ws.Range("K70").Value = 1001305.02
$ws.Range("M70").Value = -1001035.02
# row 73 (hunk 5)
$ws.Range("H73").Value = 294670.6
$ws.Range("I73").Value = 333768.34
$ws.Range("K73").Value = 1001305.02
$ws.Range("M73").Value = -1000369.02
# row 98 (hunk 6)
$ws.Range("H98").Value = 6357
$ws.Range("I98").Value = 300
$ws.Range("K98").Value = 300
$ws.Range("M98").Value = 1198
# row 122 (hunk 7)
$ws.Range("H122").Value = 6357
$ws.Range("I122").Value = 300
$ws.Range("K122").Value = 900
$ws.Range("M122").Value = 1550
# row 133 (hunk 8)
$ws.Range("H133").Value = 69999
$ws.Range("J133").Value = 69999
$ws.Range("L133").Value = 69999
$ws.Range("N133").Value = -80119
# row 137 (hunk 9)
$ws.Range("H137").Value = 406961.84
$ws.Range("I137").Value = 629318.25
$ws.Range("K137").Value = 1887954.75
$ws.Range("M137").Value = -1885404.75
# row 141 (hunk 10)
$ws.Range("H141").Value = 3194.7727
$ws.Range("I141").Value = 3289.3
$ws.Range("K141").Value = 9867.900000000001
$ws.Range("M141").Value = -4687.900000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 12 (hunk 11)
$ws.Range("H12").Value = 999.8570999999999
$ws.Range("I12").Value = 999.5
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 999.5
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -826.5
$ws.Range("N12").Value = -1346
# row 32 (hunk 12)
$ws.Range("H32").Value = 42906.426
$ws.Range("I32").Value = 37352.543
$ws.Range("J32").Value = 87337.5
$ws.Range("K32").Value = 37352.543
$ws.Range("L32").Value = 87337.5
$ws.Range("M32").Value = -37065.543
$ws.Range("N32").Value = -87911.5
# row 61 (hunk 13)
$ws.Range("H61").Value = 8624931
$ws.Range("I61").Value = 3733
$ws.Range("K61").Value = 3733
$ws.Range("M61").Value = -3521
# row 136 (hunk 14)
$ws.Range("H136").Value = 8624931
$ws.Range("I136").Value = 3733
$ws.Range("K136").Value = 11199
$ws.Range("M136").Value = -8649

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 13 (hunk 15)
$ws.Range("H13").Value = 102995
$ws.Range("J13").Value = 102995
$ws.Range("L13").Value = 102995
$ws.Range("N13").Value = -103331
# row 22 (hunk 16)
$ws.Range("H22").Value = 2878.5
$ws.Range("I22").Value = 2161.7
$ws.Range("K22").Value = 2161.7
$ws.Range("M22").Value = -1988.7
# row 94 (hunk 17)
$ws.Range("H94").Value = 602.4
$ws.Range("I94").Value = 558.2222
$ws.Range("K94").Value = 558.2222
$ws.Range("M94").Value = -107.2222
# row 99 (hunk 18)
$ws.Range("H99").Value = 5210.0557
$ws.Range("I99").Value = 6552.1
$ws.Range("J99").Value = 3532.5
$ws.Range("K99").Value = 6552.1
$ws.Range("L99").Value = 3532.5
$ws.Range("M99").Value = -5054.1
$ws.Range("N99").Value = -6528.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk 19)
$ws.Range("H31").Value = 19807.732
$ws.Range("I31").Value = 31697.8
$ws.Range("J31").Value = 13862.7
$ws.Range("K31").Value = 31697.8
$ws.Range("L31").Value = 13862.7
$ws.Range("M31").Value = -31402.8
$ws.Range("N31").Value = -14452.7
# row 34 (hunk 20)
$ws.Range("H34").Value = 19807.732
$ws.Range("I34").Value = 31697.8
$ws.Range("J34").Value = 13862.7
$ws.Range("K34").Value = 31697.8
$ws.Range("L34").Value = 13862.7
$ws.Range("M34").Value = -31495.8
$ws.Range("N34").Value = -14266.7
# row 58 (hunk 21)
$ws.Range("H58").Value = 4279.4165
$ws.Range("I58").Value = 1512.25
$ws.Range("J58").Value = 5663
$ws.Range("K58").Value = 1512.25
$ws.Range("L58").Value = 5663
$ws.Range("M58").Value = -1309.25
$ws.Range("N58").Value = -6069
# row 136 (hunk 22)
$ws.Range("H136").Value = 4279.4165
$ws.Range("I136").Value = 1512.25
$ws.Range("J136").Value = 5663
$ws.Range("K136").Value = 4536.75
$ws.Range("L136").Value = 16989
$ws.Range("M136").Value = -1986.75
$ws.Range("N136").Value = -22089
# row 141 (hunk 23)
$ws.Range("H141").Value = 331045.12
$ws.Range("J141").Value = 338952.06
$ws.Range("L141").Value = 338952.06
$ws.Range("N141").Value = -349312.06

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 113 (hunk 24)
$ws.Range("H113").Value = 2720
$ws.Range("I113").Value = 2334
$ws.Range("J113").Value = 2885.4285
$ws.Range("K113").Value = 7002
$ws.Range("L113").Value = 8656.2855
$ws.Range("M113").Value = -4832
$ws.Range("N113").Value = -12996.2855

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 80 (hunk 25)
$ws.Range("H80").Value = 5919
$ws.Range("I80").Value = 4335
$ws.Range("K80").Value = 4335
$ws.Range("M80").Value = -3337
# row 83 (hunk 26)
$ws.Range("H83").Value = 5919
$ws.Range("I83").Value = 4335
$ws.Range("K83").Value = 21675
$ws.Range("M83").Value = -16683
# row 97 (hunk 27)
$ws.Range("H97").Value = 821.44446
$ws.Range("I97").Value = 699.125
$ws.Range("K97").Value = 699.125
$ws.Range("M97").Value = -203.125
# row 113 (hunk 28)
$ws.Range("H113").Value = 4139.9443
$ws.Range("J113").Value = 5303.25
$ws.Range("L113").Value = 5303.25
$ws.Range("N113").Value = -9643.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22 (hunk 29)
$ws.Range("H22").Value = 2274.75
$ws.Range("J22").Value = 3250
$ws.Range("L22").Value = 3250
$ws.Range("N22").Value = -3840
# row 27 (hunk 30)
$ws.Range("H27").Value = 2274.75
$ws.Range("J27").Value = 3250
$ws.Range("L27").Value = 3250
$ws.Range("N27").Value = -3464
# row 82 (hunk 31)
$ws.Range("H82").Value = 1647.7037
$ws.Range("J82").Value = 1985.6428
$ws.Range("L82").Value = 1985.6428
$ws.Range("N82").Value = -2707.6428
# row 85 (hunk 32)
$ws.Range("H85").Value = 1647.7037
$ws.Range("J85").Value = 1985.6428
$ws.Range("L85").Value = 1985.6428
$ws.Range("N85").Value = -4481.6428
# row 100 (hunk 33)
$ws.Range("H100").Value = 628600.2
$ws.Range("I100").Value = 1253186.6
$ws.Range("J100").Value = 4013.75
$ws.Range("K100").Value = 1253186.6
$ws.Range("L100").Value = 4013.75
$ws.Range("M100").Value = -1252645.6
$ws.Range("N100").Value = -5095.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 96 (hunk 34)
$ws.Range("H96").Value = 2440.4285
$ws.Range("I96").Value = 2456.2
$ws.Range("J96").Value = 2401
$ws.Range("K96").Value = 2456.2
$ws.Range("L96").Value = 2401
$ws.Range("M96").Value = -1083.2
$ws.Range("N96").Value = -5147
# row 100 (hunk 35)
$ws.Range("H100").Value = 1438.7188
$ws.Range("I100").Value = 1372.1111
$ws.Range("J100").Value = 1798.4
$ws.Range("K100").Value = 2744.2222
$ws.Range("L100").Value = 3596.8
$ws.Range("M100").Value = -2203.2222
$ws.Range("N100").Value = -4678.8
# row 118 (hunk 36)
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
# row 136 (hunk 37)
$ws.Range("H136").Value = 1239990.9
$ws.Range("I136").Value = 3669.8845
$ws.Range("J136").Value = 2931798.8
$ws.Range("K136").Value = 11009.6535
$ws.Range("L136").Value = 8795396.399999999
$ws.Range("M136").Value = -8459.6535
$ws.Range("N136").Value = -8800496.399999999
